$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric-looking Price cells to Text format before assignment,
# so Excel does not auto-convert them to numbers (which would also strip
# significant trailing zeros, e.g. "0.160" -> 0.16).
$forceTextRows = @(5,6,9,10,11,12,14,18,20,21,22,23,24,25,27,31,32,33,35,37,38,39,40,41,43,45,46,47,49,51)
foreach ($r in $forceTextRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '64.203.83'
$ws.Range("E2").Value = '  +8.58%  '

# Row 3
$ws.Range("D3").Value = '3.150.99'
$ws.Range("E3").Value = '  +5.92%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").Value = '592.22'
$ws.Range("E5").Value = '  +4.70%  '

# Row 6
$ws.Range("D6").Value = '147.49'
$ws.Range("E6").Value = '  +8.07%  '

# Row 7
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").Value = '3.144.85'
$ws.Range("E8").Value = '  +5.86%  '

# Row 9
$ws.Range("D9").Value = '0.535'
$ws.Range("E9").Value = '  +3.74%  '

# Row 10
$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  +21.14%  '

# Row 11
$ws.Range("D11").Value = '5.77'
$ws.Range("E11").Value = '  +9.40%  '

# Row 12
$ws.Range("D12").Value = '0.471'
$ws.Range("E12").Value = '  +5.20%  '

# Row 13
$ws.Range("E13").Value = '  +11.93%  '

# Row 14
$ws.Range("D14").Value = '36.16'
$ws.Range("E14").Value = '  +8.02%  '

# Row 15
$ws.Range("E15").Value = '  +1.01%  '

# Row 16
$ws.Range("D16").Value = '3.673.93'
$ws.Range("E16").Value = '  +5.89%  '

# Row 17
$ws.Range("D17").Value = '64.123.01'
$ws.Range("E17").Value = '  +8.46%  '

# Row 18
$ws.Range("D18").Value = '7.18'
$ws.Range("E18").Value = '  +1.79%  '

# Row 19
$ws.Range("D19").Value = '3.150.41'
$ws.Range("E19").Value = '  +5.90%  '

# Row 20
$ws.Range("D20").Value = '474.55'
$ws.Range("E20").Value = '  +9.40%  '

# Row 21
$ws.Range("D21").Value = '14.33'
$ws.Range("E21").Value = '  +5.20%  '

# Row 22
$ws.Range("D22").Value = '0.732'
$ws.Range("E22").Value = '  +1.25%  '

# Row 23
$ws.Range("D23").Value = '7.63'
$ws.Range("E23").Value = '  +8.84%  '

# Row 24
$ws.Range("D24").Value = '13.40'
$ws.Range("E24").Value = '  +3.23%  '

# Row 25
$ws.Range("D25").Value = '82.53'
$ws.Range("E25").Value = '  +3.33%  '

# Row 26
$ws.Range("E26").Value = '  +0.12%  '

# Row 27
$ws.Range("D27").Value = '8.67'
$ws.Range("E27").Value = '  +12.63%  '

# Row 28
$ws.Range("E28").Value = '  +6.55%  '

# Row 29
$ws.Range("E29").Value = '  +0.62%  '

# Row 30
$ws.Range("E30").Value = '  +0.07%  '

# Row 31
$ws.Range("D31").Value = '6.88'
$ws.Range("E31").Value = '  +11.84%  '

# Row 32
$ws.Range("D32").Value = '27.22'
$ws.Range("E32").Value = '  +6.30%  '

# Row 33
$ws.Range("D33").Value = '0.109'
$ws.Range("E33").Value = '  +6.69%  '

# Row 34
$ws.Range("D34").Value = '0.0₃0879'
$ws.Range("E34").Value = '  +16.22%  '

# Row 35
$ws.Range("D35").Value = '2.43'
$ws.Range("E35").Value = '  +18.57%  '

# Row 36
$ws.Range("E36").Value = '  +7.36%  '

# Row 37
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").Value = '6.15'
$ws.Range("E37").Value = '  +4.91%  '

# Row 38
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").Value = '3.37'
$ws.Range("E38").Value = '  +21.14%  '

# Row 39
$ws.Range("D39").Value = '50.96'
$ws.Range("E39").Value = '  +5.18%  '

# Row 40
$ws.Range("D40").Value = '451.81'
$ws.Range("E40").Value = '  +14.79%  '

# Row 41
$ws.Range("D41").Value = '8.77'
$ws.Range("E41").Value = '  +0.86%  '

# Row 42
$ws.Range("D42").Value = '2.939.68'
$ws.Range("E42").Value = '  +8.77%  '

# Row 43
$ws.Range("D43").Value = '0.0372'
$ws.Range("E43").Value = '  +6.27%  '

# Row 44
$ws.Range("E44").Value = '  +14.11%  '

# Row 45
$ws.Range("D45").Value = '0.113'
$ws.Range("E45").Value = '  +8.06%  '

# Row 46
$ws.Range("D46").Value = '2.21'
$ws.Range("E46").Value = '  +12.49%  '

# Row 47
$ws.Range("D47").Value = '35.15'
$ws.Range("E47").Value = '  +2.50%  '

# Row 49
$ws.Range("D49").Value = '123.21'
$ws.Range("E49").Value = '  +0.59%  '

# Row 50
$ws.Range("E50").Value = '  +2.69%  '

# Row 51
$ws.Range("D51").Value = '24.95'
$ws.Range("E51").Value = '  +8.22%  '
